$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows where x_nrSteps = -2 and alienID = 36:
#   x_corrSteps -= 1, x_nrSteps = -3, alienID = 46
$rows = @(2, 9, 12, 20, 22, 28)
foreach ($r in $rows) {
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value2 = $dCell.Value2 - 1
    $ws.Cells.Item($r, 6).Value2 = -3
    $ws.Cells.Item($r, 8).Value2 = 46
}

# Update the active selection on the sheet
$ws.Range("B30").Select()

$wb.Save()
